## Add the new "262144 particles" timing block to the Input sheet.
$wb = $excel.ActiveWorkbook
$input = $wb.Worksheets.Item("Input")

$lines = @(
  "Amount of time taken by each method over the course of 1 frames, simulating 262144 particles:",
  "collideParticles() = 00493020070901 nanoseconds",
  "calculateGrav()    = 12149956944835 nanoseconds",
  "updateFocus()      = 00000010111520 nanoseconds",
  "eraseParticles()   = 00017484682538 nanoseconds",
  "updateVelAndPos()  = 00000008860862 nanoseconds",
  "drawParticles()    = 00000080240423 nanoseconds",
  "saveFrame()        = 00001615842811 nanoseconds"
)

$startRow = 154
for ($i = 0; $i -lt $lines.Length; $i++) {
  $input.Cells.Item($startRow + $i, 1).Value = $lines[$i]
}

$excel.CalculateFullRebuild()
